## The deck ships two theme palettes:
##   - "Office Theme" (the default Office colours)
##   - "Integral"      (a green/earth-tone palette)
##
## Before this edit the presentation's live design (the one driving the
## Slide Master and every slide) uses the "Integral" palette. The commit
## swaps the two palettes around so the live design instead uses the
## "Office Theme" colours.
##
## Do this through the Slide Master's Theme / ThemeColorScheme, which is
## the live, editable colour scheme for the presentation's design.

$p      = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme  = $master.Theme
$colors = $theme.ThemeColorScheme

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

# "Office Theme" palette, in msoThemeColorSchemeIndex order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink
$officeThemePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = ToRGB($officeThemePalette[$i - 1])
}

# Rename the theme / colour scheme to match the restored palette
# (no-ops on hosts that treat these names as read-only, harmless otherwise).
$theme.Name  = "Office Theme"
$colors.Name = "Office"
